# Apply the two changes captured by the commit:
#   1. Slide 5's table switches to a different table style (GUID change).
#   2. The presentation's active theme ("Integral") has its colour
#      scheme swapped for the plain "Office Theme" colour values that
#      previously only lived in the Notes Master's theme part.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{1FF89F26-268B-4072-B578-D8B86AD9A72A}")

# --- 2. Theme colour scheme swap ------------------------------------------
# The deck's Slide Master / Presentation theme (accessible through any
# slide's ThemeColorScheme) gets re-coloured from "Red Violet" (Integral)
# to the classic "Office" palette.
$colorScheme = $p.Slides.Item(1).ThemeColorScheme
$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
